$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header text "MODEL_CONDITION" -> "MODELCONDITION" wherever it
# occurs on the sheet (currently the E1 header cell).
$headerRange = $ws.Cells.Find("MODEL_CONDITION")
if ($headerRange -ne $null) {
    $headerRange.Value = "MODELCONDITION"
} else {
    $ws.Range("E1").Value = "MODELCONDITION"
}

# Delete entire column A (the leading numeric column), shifting columns
# B:F left to A:E.
$ws.Range("A1").EntireColumn.Delete()
